$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 87
$ws.Range("E2").Value = 0.8529411764705882
$ws.Range("F2").Value = 0.8529411764705882
$ws.Range("G2").Value = 0.105155020195745
$ws.Range("H2").Value = 0.08969104663754716
$ws.Range("I2").Value = 519414.7856197282
$ws.Range("J2").Value = 197093.3930108641
$ws.Range("L2").Value = 197093.3930108641
$ws.Range("M2").Value = 716508.1786305922
$ws.Range("N2").Value = 10094300.9488
$ws.Range("O2").Value = 9686560.0187
$ws.Range("P2").Value = 0.01952521467415674
$ws.Range("Q2").Value = 0.02034709872548907

# Row 3
$ws.Range("C3").Value = 102
$ws.Range("D3").Value = 85
$ws.Range("E3").Value = 0.8333333333333334
$ws.Range("F3").Value = 0.8252427184466019
$ws.Range("G3").Value = 0.1124218121923839
$ws.Range("H3").Value = 0.09277528190633627
$ws.Range("I3").Value = 589272.077115087
$ws.Range("J3").Value = 229423.6187645734
$ws.Range("L3").Value = 229423.6187645734
$ws.Range("M3").Value = 818695.6958796604
$ws.Range("N3").Value = 10505163.513564
$ws.Range("O3").Value = 10097790.355561
$ws.Range("P3").Value = 0.02183912877399267
$ws.Range("Q3").Value = 0.02272018042424762

# Row 4
$ws.Range("D4").Value = 84
$ws.Range("E4").Value = 0.8076923076923077
$ws.Range("F4").Value = 0.8076923076923077
$ws.Range("G4").Value = 0.1180246687376596
$ws.Range("H4").Value = 0.09532761705734048
$ws.Range("I4").Value = 649734.5660466086
$ws.Range("J4").Value = 252568.4906365452
$ws.Range("L4").Value = 252568.4906365452
$ws.Range("M4").Value = 902303.0566831537
$ws.Range("N4").Value = 10901225.88937092
$ws.Range("O4").Value = 10492881.53662783
$ws.Range("P4").Value = 0.02316881543412548
$ws.Range("Q4").Value = 0.02407046050743035

# Row 5
$ws.Range("G5").Value = 0.1126370358327302
$ws.Range("H5").Value = 0.09654603071376873
$ws.Range("I5").Value = 692419.1877078008
$ws.Range("J5").Value = 271336.8376955385
$ws.Range("L5").Value = 271336.8376955385
$ws.Range("M5").Value = 963756.0254033392
$ws.Range("N5").Value = 11465708.21445205
$ws.Range("O5").Value = 11055013.53112666
$ws.Range("P5").Value = 0.0236650743783563
$ws.Range("Q5").Value = 0.02454423388370882

# Row 6
$ws.Range("D6").Value = 90
$ws.Range("E6").Value = 0.8490566037735849
$ws.Range("F6").Value = 0.8490566037735849
$ws.Range("G6").Value = 0.1130249947913053
$ws.Range("H6").Value = 0.09596461821903285
$ws.Range("I6").Value = 711267.298988305
$ws.Range("J6").Value = 277699.7106510397
$ws.Range("L6").Value = 277699.7106510397
$ws.Range("M6").Value = 988967.0096393448
$ws.Range("N6").Value = 11792951.18548561
$ws.Range("O6").Value = 11378485.66166047
$ws.Range("P6").Value = 0.02354794031478937
$ws.Range("Q6").Value = 0.02440568269877443
